# #5: property aircraft done
# Update property_category values:
#  - 建物 (Building) sheet: rows 2-10, column I -> "building"
#  - 汽車 (Car) sheet:      rows 2-3,  column H -> "car"

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
for ($r = 2; $r -le 10; $r++) {
    $wsBuilding.Cells.Item($r, 9).Value = "building"
}

$wsCar = $wb.Worksheets.Item("汽車")
for ($r = 2; $r -le 3; $r++) {
    $wsCar.Cells.Item($r, 8).Value = "car"
}
